$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 7 on every sheet corresponds to file 910bdb67-0a00-40aa-bd98-949a431a4405.md
# which is ready for handoff now (status changed from "In Translation" to "Ready for handoff"),
# together with refreshed "Latest Handoff"/"Latest HO Xliff Generate" timestamps.

$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-08-26 16:45:26"

$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("H7").Value = "2016-08-26 16:45:22"

$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("H7").Value = "2016-08-26 16:45:26"
